# Add standard deviation (SD) values to table 1 numeric cells.
# Each original mean value becomes "<new mean> +/- <SD>" (with a couple of
# cells also gaining a trailing ')' per the source data).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 = Industrial
$t.Cell(2,2).Range.Text = "7683 +/- 2576"
$t.Cell(2,3).Range.Text = "-1849 +/- 2440)"
$t.Cell(2,4).Range.Text = "-7.3 +/- 7.5"
$t.Cell(2,5).Range.Text = "-10.4 +/- 10.4"
$t.Cell(2,6).Range.Text = "3.5 +/- 2"
$t.Cell(2,7).Range.Text = "6.8 +/- 13.6"

# Row 3 = Agricultural
$t.Cell(3,2).Range.Text = "2537 +/- 6041"
$t.Cell(3,3).Range.Text = "1230 +/- 3855)"
$t.Cell(3,4).Range.Text = "28.3 +/- 39.6"
$t.Cell(3,5).Range.Text = "-8.6 +/- 9.9"
$t.Cell(3,6).Range.Text = "2.1 +/- 2.9"
$t.Cell(3,7).Range.Text = "14.4 +/- 11"

# Row 4 = Forested
$t.Cell(4,2).Range.Text = "1571 +/- 3498"
$t.Cell(4,3).Range.Text = "698 +/- 2908)"
$t.Cell(4,4).Range.Text = "21.2 +/- 28.6"
$t.Cell(4,5).Range.Text = "-3.4 +/- 11"
$t.Cell(4,6).Range.Text = "3 +/- 3.1"
$t.Cell(4,7).Range.Text = "5.1 +/- 12.2"
